$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.503.26'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.87%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.402.91'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.45%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.07%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.68'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.52%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.401.12'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.45%  '

$ws.Range("E9").Value = '  -2.91%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.91'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.40%  '

$ws.Range("E11").Value = '  -6.50%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.403'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.974.66'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.51%  '

$ws.Range("E14").Value = '  -7.12%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.34'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -7.20%  '

$ws.Range("E16").Value = '  -0.64%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.400.92'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.24%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '65.548.52'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.91%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.32'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.73%  '

$ws.Range("E20").Value = '  -5.99%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.53'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.73%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '413.20'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.576'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.83%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '76.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.43%  '

$ws.Range("E25").Value = '  +0.11%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.540.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.31%  '

$ws.Range("E27").Value = '  -9.84%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.17'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.43%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.75'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.80%  '

$ws.Range("E30").Value = '  -3.41%  '

$ws.Range("E31").Value = '  -0.05%  '

$ws.Range("E32").Value = '  -5.75%  '

$ws.Range("E33").Value = '  -9.18%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.37'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.25%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.398.21'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.28%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.50'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -9.20%  '

$ws.Range("E38").Value = '  -7.55%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.00%  '

$ws.Range("E40").Value = '  -5.91%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '167.90'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0848'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.79%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.869'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.78%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.01'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.90'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -10.94%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '45.32'
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.31'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -9.46%  '

$ws.Range("E48").Value = '  -4.85%  '

$ws.Range("E49").Value = '  -6.19%  '

$ws.Range("E50").Value = '  -7.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.913'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.27%  '
